$d = $word.ActiveDocument

# The four list items that got "cancelled" (struck through) in this edit:
#   - "register del consumer"
#   - "login del consumer"
#   - "register admin"
#   - "login admin"
# Apply strikethrough formatting to the whole paragraph (including the
# paragraph mark) so both the run(s) and the paragraph mark's rPr pick up
# <w:strike/>, matching how Word applies character formatting to a
# selected paragraph.

$targets = @("register del consumer", "login del consumer", "register admin", "login admin")

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($targets -contains $text) {
        $p.Range.Font.StrikeThrough = 1
    }
}
